# Fruta / hortaliza, semanal
#
# Insert two new weekly price records (row 78 and 79 in the final sheet)
# for "Santina" cherries at Terminal Hortofruticola Agro Chillan, pushing
# the previously existing rows 78-100 down to 80-102 (dimension grows from
# A1:T100 to A1:T102).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the two new rows right above the current row 78.
$ws.Rows("78:79").Insert()

# New row 78: Santina / Primera
$ws.Range("A78").Value = 7
$ws.Range("B78").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C78").Value = "Ñuble"
$ws.Range("D78").Value = 44588
$ws.Range("E78").Value = 16
$ws.Range("F78").Value = "Fruta"
$ws.Range("G78").Value = 100103
$ws.Range("H78").Value = "Frutos de hueso (carozo)"
$ws.Range("I78").Value = 100103001
$ws.Range("J78").Value = "Cereza"
$ws.Range("K78").Value = "Santina"
$ws.Range("L78").Value = "Primera"
$ws.Range("M78").Value = 160
$ws.Range("N78").Value = 5500
$ws.Range("O78").Value = 6000
$ws.Range("P78").Value = 5750
$ws.Range("Q78").Value = "$/bandeja 10 kilos"
$ws.Range("R78").Value = "Provincia de Curicó"
$ws.Range("S78").Value = 575
$ws.Range("T78").Value = 10

# New row 79: Santina / Segunda
$ws.Range("A79").Value = 7
$ws.Range("B79").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C79").Value = "Ñuble"
$ws.Range("D79").Value = 44588
$ws.Range("E79").Value = 16
$ws.Range("F79").Value = "Fruta"
$ws.Range("G79").Value = 100103
$ws.Range("H79").Value = "Frutos de hueso (carozo)"
$ws.Range("I79").Value = 100103001
$ws.Range("J79").Value = "Cereza"
$ws.Range("K79").Value = "Santina"
$ws.Range("L79").Value = "Segunda"
$ws.Range("M79").Value = 120
$ws.Range("N79").Value = 4500
$ws.Range("O79").Value = 5000
$ws.Range("P79").Value = 4750
$ws.Range("Q79").Value = "$/bandeja 10 kilos"
$ws.Range("R79").Value = "Provincia de Curicó"
$ws.Range("S79").Value = 475
$ws.Range("T79").Value = 10
